# ---------------------------------------------------------------
# feat: add 2022-Q1 data
#
# 1. The sheet that used to be named "总计" becomes the new "2022-Q1"
#    fund-ranking sheet (keeps its sheetId so the workbook.xml sheetId
#    numbering matches: 2022-Q1 -> 4).
# 2. A brand-new sheet is created and named "总计" (gets the next
#    sheetId = 5) and moved to the last tab position.
# 3. "总计" is repopulated with its original 4 rows, plus a new first
#    data row summarizing 2022-Q1.
# ---------------------------------------------------------------

$wb = $excel.ActiveWorkbook

$fundData = @(
    @('510500', '南方中证500ETF', '402.81', '95.26', '0.57', '2.2960', 5),
    @('515220', '国泰中证煤炭ETF', '32.88', '99.67', '6.71', '2.2062', 5),
    @('161032', '富国中证煤炭指数', '21.64', '93.99', '6.96', '1.5061', 4),
    @('168204', '中融中证煤炭指数', '8.68', '92.62', '6.83', '0.5928', 4),
    @('512500', '华夏中证500ETF', '60.13', '98.15', '0.60', '0.3608', 5),
    @('159922', '嘉实中证500ETF', '36.91', '97.34', '0.59', '0.2178', 5),
    @('510510', '广发中证500ETF', '27.79', '98.59', '0.60', '0.1667', 6),
    @('159820', '天弘中证500ETF', '23.75', '95.88', '0.58', '0.1378', 5),
    @('159930', '汇添富中证能源ETF', '2.12', '99.19', '6.15', '0.1304', 6),
    @('510580', '易方达中证500ETF', '16.77', '95.49', '0.58', '0.0973', 6),
    @('510590', '平安中证500ETF', '12.98', '98.06', '0.60', '0.0779', 5),
    @('009613', '上银中证500指数增强A', '2.83', '90.41', '1.43', '0.0405', 2),
    @('159968', '博时中证500ETF', '6.22', '95.76', '0.58', '0.0361', 6),
    @('512510', '华泰柏瑞中证500ETF', '5.50', '96.08', '0.58', '0.0319', 5),
    @('160616', '鹏华中证500指数(LOF) A', '5.45', '94.90', '0.58', '0.0316', 5),
    @('159945', '广发中证全指能源ETF', '0.53', '98.90', '4.68', '0.0248', 6),
    @('009614', '上银中证500指数增强C', '1.70', '90.41', '1.43', '0.0243', 2),
    @('159982', '鹏华中证500ETF', '3.91', '94.84', '0.58', '0.0227', 5),
    @('510530', '工银瑞信中证500ETF', '3.64', '99.00', '0.60', '0.0218', 6),
    @('165511', '信诚中证500指数（LOF）A', '2.78', '93.31', '0.57', '0.0158', 4),
    @('510560', '国寿安保中证500ETF', '2.42', '99.37', '0.61', '0.0148', 5),
    @('010992', '西藏东财中证500指数A', '2.48', '94.93', '0.58', '0.0144', 5),
    @('561350', '国泰中证500ETF', '1.78', '96.89', '0.58', '0.0103', 9),
    @('510440', '大成中证500沪市ETF', '0.41', '96.76', '1.10', '0.0045', 3),
    @('159935', '景顺长城中证500ETF', '0.74', '98.01', '0.60', '0.0044', 5),
    @('010993', '西藏东财中证500指数C', '0.72', '94.93', '0.58', '0.0042', 5),
    @('515190', '中银证券中证500ETF', '0.62', '99.21', '0.61', '0.0038', 6),
    @('660011', '农银中证500指数', '0.55', '94.34', '0.58', '0.0032', 5),
    @('006611', '人保中证500指数', '0.44', '92.48', '0.69', '0.0030', 4),
    @('159999', '永赢中证500ETF', '0.49', '97.41', '0.59', '0.0029', 6),
    @('007943', '富安达中证 500 指数增强', '0.21', '93.50', '1.38', '0.0029', 6),
    @('001351', '诺安中证500指数增强A', '0.48', '94.42', '0.56', '0.0027', 6),
    @('515530', '泰康中证500ETF', '0.47', '94.55', '0.58', '0.0027', 5),
    @('510550', '方正富邦中证500ETF', '0.24', '97.63', '0.60', '0.0014', 6),
    @('006938', '鹏华中证500指数(LOF) C', '0.25', '94.90', '0.58', '0.0014', 5),
    @('515550', '中融中证500ETF', '0.23', '91.02', '0.57', '0.0013', 5),
    @('013119', '信诚中证500指数（LOF）C', '0.12', '93.31', '0.57', '0.0007', 4),
    @('510570', '兴业中证500ETF', '0.10', '96.12', '0.58', '0.0006', 8),
    @('010355', '诺安中证500指数增强C', '0.04', '94.42', '0.56', '0.0002', 6)
)


$totalsData = @(
    @('2022-Q1', 39, 8.119999999999999),
    @('2021-Q4', 37, 10.8),
    @('2021-Q3', 1, 0.02),
    @('2021-Q2', 1, 0.01)
)

# A cell elsewhere in the workbook that already carries the bold +
# bordered + centered "header/index" style, used purely as a format
# donor via PasteSpecial so the new sheets pick up the identical style.
$styleDonor = $wb.Worksheets.Item("2021-Q4").Range("B1")

# -----------------------------------------------------------------
# Step 1: turn the old "总计" sheet into the new "2022-Q1" sheet.
# -----------------------------------------------------------------
$fundSheet = $wb.Worksheets.Item("总计")
$fundSheet.UsedRange.Clear()
$fundSheet.Name = "2022-Q1"

$fundHeaders = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($col = 0; $col -lt $fundHeaders.Length; $col++) {
    $cell = $fundSheet.Cells.Item(1, $col + 2)
    $cell.Value = $fundHeaders[$col]
}
$headerRange = $fundSheet.Range("B1:H1")
$styleDonor.Copy()
$headerRange.PasteSpecial(-4122)

# Columns B..F and G hold text-like numbers (fund codes / percentages /
# scale) that must be stored as text, not auto-converted numbers -
# format the block as Text before writing.
$fundSheet.Range("B2:G40").NumberFormat = "@"

for ($i = 0; $i -lt $fundData.Length; $i++) {
    $row = $i + 2
    $entry = $fundData[$i]
    $idxCell = $fundSheet.Cells.Item($row, 1)
    $idxCell.Value = $i
    $styleDonor.Copy()
    $idxCell.PasteSpecial(-4122)

    $fundSheet.Cells.Item($row, 2).Value = $entry[0]
    $fundSheet.Cells.Item($row, 3).Value = $entry[1]
    $fundSheet.Cells.Item($row, 4).Value = $entry[2]
    $fundSheet.Cells.Item($row, 5).Value = $entry[3]
    $fundSheet.Cells.Item($row, 6).Value = $entry[4]
    $fundSheet.Cells.Item($row, 7).Value = $entry[5]
    $fundSheet.Cells.Item($row, 8).Value = $entry[6]
}

# -----------------------------------------------------------------
# Step 2: create a fresh "总计" sheet after "2022-Q1" and fill it in.
# -----------------------------------------------------------------
$totalSheet = $wb.Worksheets.Add()
$totalSheet.Name = "总计"
$totalSheet.Move($null, $wb.Worksheets.Item($wb.Worksheets.Count))

$totalSheet.Cells.Item(1, 2).Value = "日期"
$totalSheet.Cells.Item(1, 3).Value = "持有数量(只)"
$totalSheet.Cells.Item(1, 4).Value = "持有市值(亿元)"
$totalHeaderRange = $totalSheet.Range("B1:D1")
$styleDonor.Copy()
$totalHeaderRange.PasteSpecial(-4122)

for ($i = 0; $i -lt $totalsData.Length; $i++) {
    $row = $i + 2
    $entry = $totalsData[$i]
    $idxCell = $totalSheet.Cells.Item($row, 1)
    $idxCell.Value = $i
    $styleDonor.Copy()
    $idxCell.PasteSpecial(-4122)

    $totalSheet.Cells.Item($row, 2).Value = $entry[0]
    $totalSheet.Cells.Item($row, 3).Value = $entry[1]
    $totalSheet.Cells.Item($row, 4).Value = $entry[2]
}
